# Generate Report for Handoff
#
# Marks file "b.md" as ready for handoff in the localization status report:
#   - Overview sheet: status + latest-handoff columns for the "b.md" row
#   - zh-cn sheet: status, latest handoff file, and latest handoff datetime
#     for the "b.md" row (new handoff package
#     b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf)
#   - de-de sheet: same, for the de-de handoff package

$wb = $excel.ActiveWorkbook

$status = "Ready for handoff"

# --- Overview sheet -------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $status
$overview.Range("C3").Value = $status
$overview.Range("D3").Value = "2016-19-11 08:19:59"

# --- zh-cn sheet ------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $status
$zhcn.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("E3").Value = "2016-03-11 08:19:55"

# --- de-de sheet ------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $status
$dede.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("E3").Value = "2016-03-11 08:19:59"
